# Apply BOM update: fill in missing positions for existing parts, and add
# rows for the new LED-matrix related components (diode, decoupling caps,
# ESP32 module, current-limit resistor, additional bypass cap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows -------------------------------------------------

# Row 3 (C 4,7u) was missing its "Positionen" (designator) value.
$ws.Range("D3").Value = "C1, C2"

# Row 4 (D / RGB LEDs) had a placeholder "D1 - D42" Positionen value that is
# no longer accurate now that other diodes (D1) exist separately - clear it.
$ws.Range("D4").Value = ""

# --- New rows --------------------------------------------------------------

# Row 6: Schottky diode D1
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "D"
$ws.Range("C6").Value = "40V 3A Schottky"
$ws.Range("D6").Value = "D1"
$ws.Range("E6").Value = "B340A-13-F DII"

# Row 7: 470u bulk capacitors C3, C5
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "C"
$ws.Range("C7").Value = "470u"
$ws.Range("D7").Value = "C3, C5"
$ws.Range("E7").Value = "FK 470/6,3 SP"

# Row 8: 100p capacitors C4, C6
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "100p"
$ws.Range("D8").Value = "C4, C6"
$ws.Range("E8").Value = "FK 470/6,3 SP"

# Row 9: ESP32 microcontroller
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = " µC"
$ws.Range("C9").Value = "ESP32 WROOM 32D"
$ws.Range("D9").Value = "U1"
$ws.Range("F9").Value = "1925467 - VQ"

# Row 10: current-limiting resistors
$ws.Range("A10").Value = 42
$ws.Range("B10").Value = "R"
$ws.Range("C10").Value = "91R"
$ws.Range("E10").Value = "RND 155HP05 EQ"

# Row 11: 100n bypass capacitors
$ws.Range("A11").Value = 42
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "100n"
$ws.Range("E11").Value = "KEM X7R0805 100N"

# --- Cosmetic updates matching the saved workbook state --------------------

# Column C widened slightly (now holding longer values like "ESP32 WROOM 32D").
$ws.Columns.Item(3).ColumnWidth = 17.14

# Put the selection where the user last left it after entering the new data.
$ws.Range("E11").Select()
